$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("Y25").Value = 64
$ws.Range("AB25").Value = 464
$ws.Range("Y26").Value = 74
$ws.Range("AB26").Value = 522
$ws.Range("Y27").Value = 81
$ws.Range("AB27").Value = 565
$ws.Range("Y28").Value = 90
$ws.Range("AB28").Value = 634
$ws.Range("Y29").Value = 94
$ws.Range("AB29").Value = 667
$ws.Range("Y30").Value = 95
$ws.Range("AB30").Value = 702
$ws.Range("AB31").Value = 778
$ws.Range("AB32").Value = 851
$ws.Range("AB33").Value = 981
$ws.Range("AB34").Value = 1093
$ws.Range("AB35").Value = 1183
$ws.Range("AB36").Value = 1212
$ws.Range("K37").Value = 11
$ws.Range("AB37").Value = 1247
$ws.Range("AB38").Value = 1276

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("Y29").Value = 1
$ws.Range("Y30").Value = 1
$ws.Range("K37").Value = 1

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("Y6").Value = -1
$ws.Range("AB6").Value = -28
$ws.Range("Y7").Value = -1
$ws.Range("AB7").Value = -31
$ws.Range("Y8").Value = -1
$ws.Range("AB8").Value = -40
$ws.Range("Y9").Value = -1
$ws.Range("AB9").Value = -42
$ws.Range("Y10").Value = -1
$ws.Range("AB10").Value = -56
$ws.Range("Y11").Value = -2
$ws.Range("AB11").Value = -55
$ws.Range("Y12").Value = -3
$ws.Range("AB12").Value = -58
$ws.Range("Y13").Value = -3
$ws.Range("AB13").Value = -67
$ws.Range("Y14").Value = -3
$ws.Range("AB14").Value = -61
$ws.Range("Y15").Value = -5
$ws.Range("AB15").Value = -65
$ws.Range("Y16").Value = -5
$ws.Range("AB16").Value = -65
$ws.Range("Y17").Value = -6
$ws.Range("AB17").Value = -77
$ws.Range("Y18").Value = -7
$ws.Range("AB18").Value = -81
$ws.Range("Y19").Value = -7
$ws.Range("AB19").Value = -77
$ws.Range("Y20").Value = -8
$ws.Range("AB20").Value = -76
$ws.Range("Y21").Value = -8
$ws.Range("AB21").Value = -80
$ws.Range("Y22").Value = -7
$ws.Range("AB22").Value = -80
$ws.Range("Y23").Value = -6
$ws.Range("AB23").Value = -78
$ws.Range("Y24").Value = -7
$ws.Range("AB24").Value = -73
$ws.Range("Y25").Value = -8
$ws.Range("AB25").Value = -79
$ws.Range("Y26").Value = -8
$ws.Range("AB26").Value = -79
$ws.Range("Y27").Value = -10
$ws.Range("AB27").Value = -82
$ws.Range("Y28").Value = -10
$ws.Range("AB28").Value = -80
$ws.Range("Y29").Value = -11
$ws.Range("AB29").Value = -81
$ws.Range("Y30").Value = -11
$ws.Range("AB30").Value = -82
$ws.Range("AB31").Value = -73
$ws.Range("AB32").Value = -71
$ws.Range("AB33").Value = -71
$ws.Range("AB34").Value = -77
$ws.Range("AB35").Value = -78
$ws.Range("AB36").Value = -76
$ws.Range("K37").Value = -2
$ws.Range("AB37").Value = -73
$ws.Range("AB38").Value = -75

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("Y29").Value = -4
$ws.Range("Y30").Value = -4
$ws.Range("K37").Value = -1

$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("Y29").Value = -2
$ws.Range("Y30").Value = -2
$ws.Range("K37").Value = -1

$ws = $wb.Worksheets.Item("Released")
$ws.Range("Y6").Value = 3
$ws.Range("AB6").Value = 67
$ws.Range("Y7").Value = 3
$ws.Range("AB7").Value = 71
$ws.Range("Y8").Value = 3
$ws.Range("AB8").Value = 75
$ws.Range("Y9").Value = 3
$ws.Range("AB9").Value = 77
$ws.Range("Y10").Value = 3
$ws.Range("AB10").Value = 85
$ws.Range("Y11").Value = 4
$ws.Range("AB11").Value = 93
$ws.Range("Y12").Value = 5
$ws.Range("AB12").Value = 102
$ws.Range("Y13").Value = 5
$ws.Range("AB13").Value = 113
$ws.Range("Y14").Value = 5
$ws.Range("AB14").Value = 115
$ws.Range("Y15").Value = 7
$ws.Range("AB15").Value = 119
$ws.Range("Y16").Value = 7
$ws.Range("AB16").Value = 121
$ws.Range("Y17").Value = 8
$ws.Range("AB17").Value = 134
$ws.Range("Y18").Value = 9
$ws.Range("AB18").Value = 138
$ws.Range("Y19").Value = 10
$ws.Range("AB19").Value = 141
$ws.Range("Y20").Value = 11
$ws.Range("AB20").Value = 160
$ws.Range("Y21").Value = 12
$ws.Range("AB21").Value = 163
$ws.Range("Y22").Value = 12
$ws.Range("AB22").Value = 163
$ws.Range("Y23").Value = 12
$ws.Range("AB23").Value = 164
$ws.Range("Y24").Value = 13
$ws.Range("AB24").Value = 166
$ws.Range("Y25").Value = 14
$ws.Range("AB25").Value = 168
$ws.Range("Y26").Value = 15
$ws.Range("AB26").Value = 178
$ws.Range("Y27").Value = 17
$ws.Range("AB27").Value = 181
$ws.Range("Y28").Value = 17
$ws.Range("AB28").Value = 183
$ws.Range("Y29").Value = 18
$ws.Range("AB29").Value = 186
$ws.Range("Y30").Value = 18
$ws.Range("AB30").Value = 187
$ws.Range("AB31").Value = 188
$ws.Range("AB32").Value = 192
$ws.Range("AB33").Value = 200
$ws.Range("AB34").Value = 204
$ws.Range("AB35").Value = 204
$ws.Range("AB36").Value = 204
$ws.Range("AB37").Value = 204
$ws.Range("AB38").Value = 204
